$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 823.36365
$ws.Range("I12").Value = 794.8889
$ws.Range("K12").Value = 794.8889
$ws.Range("M12").Value = -624.8889

$ws.Range("H42").Value = 33.75
$ws.Range("I42").Value = 31.666666
$ws.Range("J42").Value = 40
$ws.Range("K42").Value = 94.99999800000001
$ws.Range("L42").Value = 120
$ws.Range("M42").Value = 135.000002
$ws.Range("N42").Value = -580

$ws.Range("H43").Value = 4700.5
$ws.Range("I43").Value = 4500
$ws.Range("J43").Value = 4901
$ws.Range("K43").Value = 4500
$ws.Range("L43").Value = 4901
$ws.Range("M43").Value = -4431
$ws.Range("N43").Value = -5039

$ws.Range("H116").Value = 1500
$ws.Range("I116").Value = 1500
$ws.Range("K116").Value = 1500
$ws.Range("M116").Value = 1942

$ws.Range("H137").Value = 2397.8
$ws.Range("I137").Value = 1152.75
$ws.Range("J137").Value = 3820.7144
$ws.Range("K137").Value = 3458.25
$ws.Range("L137").Value = 11462.1432
$ws.Range("M137").Value = -908.25
$ws.Range("N137").Value = -16562.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2985.9092
$ws.Range("I35").Value = 2534.5
$ws.Range("J35").Value = 7500
$ws.Range("K35").Value = 2534.5
$ws.Range("L35").Value = 7500
$ws.Range("M35").Value = -2128.5
$ws.Range("N35").Value = -8312

$ws.Range("H44").Value = 11572.05
$ws.Range("J44").Value = 11865.315
$ws.Range("L44").Value = 11865.315
$ws.Range("N44").Value = -12841.315

$ws.Range("H122").Value = 2034.8572
$ws.Range("I122").Value = 1311.75
$ws.Range("K122").Value = 3935.25
$ws.Range("M122").Value = -1485.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 660.3333
$ws.Range("J22").Value = 485
$ws.Range("L22").Value = 485
$ws.Range("N22").Value = -831

$ws.Range("H26").Value = 39893.8
$ws.Range("I26").Value = 39893.8
$ws.Range("K26").Value = 39893.8
$ws.Range("M26").Value = -39601.8

$ws.Range("H96").Value = 24749.5
$ws.Range("I96").Value = 24749.5
$ws.Range("K96").Value = 24749.5
$ws.Range("M96").Value = -22003.5

$ws.Range("H99").Value = 200001310
$ws.Range("I99").Value = 200001310
$ws.Range("K99").Value = 200001310
$ws.Range("M99").Value = -199999812

$ws.Range("H105").Value = 4330749
$ws.Range("I105").Value = 6495087.5
$ws.Range("K105").Value = 6495087.5
$ws.Range("M105").Value = -6493340.5

$ws.Range("H107").Value = 50007300
$ws.Range("J107").Value = 8876.625
$ws.Range("L107").Value = 8876.625
$ws.Range("N107").Value = -12716.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 3000
$ws.Range("J8").Value = 3000
$ws.Range("L8").Value = 3000
$ws.Range("N8").Value = -3280

$ws.Range("H16").Value = 2410
$ws.Range("I16").Value = 2410
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2410
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2123
$ws.Range("N16").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H88").Value = 6909.8
$ws.Range("J88").Value = 6909.8
$ws.Range("L88").Value = 6909.8
$ws.Range("N88").Value = -7721.8

$ws.Range("H91").Value = 6909.8
$ws.Range("J91").Value = 6909.8
$ws.Range("L91").Value = 6909.8
$ws.Range("N91").Value = -9717.799999999999

$ws.Range("H105").Value = 2682.7778
$ws.Range("I105").Value = 1189
$ws.Range("K105").Value = 1189
$ws.Range("M105").Value = 558

$ws.Range("H107").Value = 1325.2941
$ws.Range("I107").Value = 545.7
$ws.Range("K107").Value = 545.7
$ws.Range("M107").Value = 1374.3

$ws.Range("H113").Value = 2410
$ws.Range("I113").Value = 2410
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2410
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -240
$ws.Range("N113").ClearContents()

$ws.Range("H134").Value = 2677.4
$ws.Range("I134").Value = 2677.4
$ws.Range("K134").Value = 8032.200000000001
$ws.Range("M134").Value = -5497.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 8986.666999999999
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 8986.666999999999
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 26960.001
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -27548.001

$ws.Range("H58").Value = 950
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 989.6
$ws.Range("I31").Value = 989.6
$ws.Range("K31").Value = 989.6
$ws.Range("M31").Value = -697.6

$ws.Range("H37").Value = 989.6
$ws.Range("I37").Value = 989.6
$ws.Range("K37").Value = 989.6
$ws.Range("M37").Value = -712.6

$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H63").Value = 31199.8
$ws.Range("J63").Value = 31199.8
$ws.Range("L63").Value = 31199.8
$ws.Range("N63").Value = -32571.8

$ws.Range("H66").Value = 31199.8
$ws.Range("J66").Value = 31199.8
$ws.Range("L66").Value = 93599.39999999999
$ws.Range("N66").Value = -100463.4

$ws.Range("H80").Value = 4663
$ws.Range("I80").Value = 4495
$ws.Range("K80").Value = 4495
$ws.Range("M80").Value = -3497

$ws.Range("H83").Value = 4663
$ws.Range("I83").Value = 4495
$ws.Range("K83").Value = 22475
$ws.Range("M83").Value = -17483

$ws.Range("H102").Value = 4224.4165
$ws.Range("I102").Value = 3615.6667
$ws.Range("K102").Value = 3615.6667
$ws.Range("M102").Value = -1993.6667

$ws.Range("H107").Value = 1200.5
$ws.Range("I107").Value = 1200.5
$ws.Range("K107").Value = 1200.5
$ws.Range("M107").Value = 719.5

$ws.Range("H122").Value = 2912.6
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 3265.75
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 9797.25
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -14697.25

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 7078.5557
$ws.Range("I132").Value = 7078.5557
$ws.Range("K132").Value = 21235.6671
$ws.Range("M132").Value = -18705.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1301.4286
$ws.Range("I7").Value = 1302.2
$ws.Range("J7").Value = 1299.5
$ws.Range("K7").Value = 1302.2
$ws.Range("L7").Value = 1299.5
$ws.Range("M7").Value = -1190.2
$ws.Range("N7").Value = -1523.5

$ws.Range("H22").Value = 991.75
$ws.Range("J22").Value = 1055.6666
$ws.Range("L22").Value = 1055.6666
$ws.Range("N22").Value = -1645.6666

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()

$ws.Range("H27").Value = 991.75
$ws.Range("J27").Value = 1055.6666
$ws.Range("L27").Value = 1055.6666
$ws.Range("N27").Value = -1269.6666

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H100").Value = 7447.3335
$ws.Range("I100").Value = 4256.5
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 4256.5
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -3715.5
$ws.Range("N100").Value = -11082

$ws.Range("H126").Value = 1301.4286
$ws.Range("I126").Value = 1302.2
$ws.Range("J126").Value = 1299.5
$ws.Range("K126").Value = 3906.6
$ws.Range("L126").Value = 3898.5
$ws.Range("M126").Value = -1436.6
$ws.Range("N126").Value = -8838.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 950
$ws.Range("J22").Value = 950
$ws.Range("L22").Value = 950
$ws.Range("N22").Value = -1536

$ws.Range("H62").Value = 7309.875
$ws.Range("I62").Value = 2163.3333
$ws.Range("J62").Value = 10397.8
$ws.Range("K62").Value = 2163.3333
$ws.Range("L62").Value = 10397.8
$ws.Range("M62").Value = -1539.3333
$ws.Range("N62").Value = -11645.8

$ws.Range("H65").Value = 7309.875
$ws.Range("I65").Value = 2163.3333
$ws.Range("J65").Value = 10397.8
$ws.Range("K65").Value = 10816.6665
$ws.Range("L65").Value = 51989
$ws.Range("M65").Value = -7696.666499999999
$ws.Range("N65").Value = -58229

$ws.Range("H107").Value = 27778316
$ws.Range("I107").Value = 27778316
$ws.Range("K107").Value = 83334948
$ws.Range("M107").Value = -83333028

$ws.Range("H122").Value = 2470
$ws.Range("I122").Value = 2289.2856
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6867.8568
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4417.8568
$ws.Range("N122").Value = -19900

$ws.Range("H126").Value = 3566.0417
$ws.Range("J126").Value = 5668
$ws.Range("L126").Value = 17004
$ws.Range("N126").Value = -21944

$ws.Range("H132").Value = 1940.3478
$ws.Range("I132").Value = 1839.4286
$ws.Range("K132").Value = 5518.2858
$ws.Range("M132").Value = -2988.2858
